## Adds the "route set up and flashcard dapp" notes paragraphs.
## The empty paragraph that used to sit right after the "...component did
## mount." sentence is replaced by three fully-populated paragraphs.

$d = $word.ActiveDocument

# Locate the end of the sentence that precedes the empty paragraph we need
# to fill in, then move to the paragraph right after it.
$rng = $d.Content
$found = $rng.Find.Execute(
    "component did mount.",  # Find What
    $true,                   # MatchCase
    $false,                  # MatchWholeWord
    $false,                  # MatchWildcards
    $false,                  # MatchSoundsLike
    $false,                  # MatchAllWordForms
    $true,                   # Forward
    1,                       # Wrap (wdFindContinue)
    $false,                  # Format
    "",                      # ReplaceWith
    0                        # Replace (wdReplaceNone)
)

if (-not $found) {
    throw "Could not find anchor sentence for the insertion point."
}

$rng.Collapse(0)            # wdCollapseEnd
$anchorPara = $rng.Paragraphs(1)
$targetPara = $anchorPara.Next()

# Sanity check: the paragraph we are about to replace should be empty.
if ($targetPara.Range.Text.Trim().Length -ne 0) {
    throw "Target paragraph was not empty; aborting to avoid clobbering content."
}

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Includes is an es6 method and checks if a string is in another string. </w:t></w:r><w:r><w:t xml:space="preserve">It returns a Bool. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Remmeber</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> we can do react </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fontaweosme</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">React Router is a collection of navigational components. </w:t></w:r><w:r><w:t xml:space="preserve">Can help create </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bookmarable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> URLS or deep links. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NavLink</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> just gives us a lot more styling power. </w:t></w:r><w:r><w:t xml:space="preserve">We can call render in Route and pass a callback. </w:t></w:r><w:r><w:t xml:space="preserve">We can render a component and use this in the Route </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ehwne</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> we need to pass props into the component. </w:t></w:r></w:p>
'@

# InsertXML replaces the contents of the exact range it targets (here the
# whole empty paragraph, paragraph mark included), so this swaps the single
# empty <w:p/> for the three new paragraphs described in the diff.
[void]$targetPara.Range.InsertXML($newParagraphsXml)
